# "nuts2 rest sector basically done"
# Adds a new "Brennstoff allgemein" (general fuel) row to the Data and
# Data_final sheets, and moves the active tab/selection over to Data_final
# (row 13, column C) as the user's next point of work, leaving a new
# selection of A12:E12 behind on the Data sheet.

$wb = $excel.ActiveWorkbook

$wsData  = $wb.Worksheets.Item("Data")
$wsFinal = $wb.Worksheets.Item("Data_final")

# --- Data sheet: new row 12 ("Brennstoff allgemein") -----------------------
# Pull formatting from row 8 (same shape: A/C/D/E populated, B empty) one
# column at a time so we don't stamp an empty styled cell into B12.
$wsData.Range("A8").Copy() | Out-Null
$wsData.Range("A12").PasteSpecial(-4122) | Out-Null
$wsData.Range("C8").Copy() | Out-Null
$wsData.Range("C12").PasteSpecial(-4122) | Out-Null
$wsData.Range("D8").Copy() | Out-Null
$wsData.Range("D12").PasteSpecial(-4122) | Out-Null
$wsData.Range("E8").Copy() | Out-Null
$wsData.Range("E12").PasteSpecial(-4122) | Out-Null

$wsData.Range("A12").Value = "Brennstoff allgemein"
$wsData.Range("C12").Value = 0.9
$wsData.Range("D12").Value = 0
$wsData.Range("E12").Value = "own assumption"

# --- Data_final sheet: new row 12 ("Brennstoff allgemein") -----------------
$wsFinal.Range("A8").Copy() | Out-Null
$wsFinal.Range("A12").PasteSpecial(-4122) | Out-Null
$wsFinal.Range("C8").Copy() | Out-Null
$wsFinal.Range("C12").PasteSpecial(-4122) | Out-Null
$wsFinal.Range("D8").Copy() | Out-Null
$wsFinal.Range("D12").PasteSpecial(-4122) | Out-Null

$wsFinal.Range("A12").Value = "Brennstoff allgemein"
$wsFinal.Range("C12").Value = 1
$wsFinal.Range("D12").Value = 0

# --- Selections / active sheet ---------------------------------------------
$wsData.Range("A12:E12").Select() | Out-Null
$wsFinal.Range("C13").Select() | Out-Null

# Data_final becomes the active (front-most) tab.
$wsFinal.Activate() | Out-Null
